# Updated symbol list on Sun Jan  1 07:29:07 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (avoids Excel auto-converting
# numeric-looking / percent-looking strings into numbers).
function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Row 2 (BNB) ---
Set-TextCell "D2" "244.01"
Set-TextCell "E2" "-0.48%"

# --- Row 3 (OKB) ---
Set-TextCell "D3" "26.37"
Set-TextCell "E3" "3.10%"

# --- Row 4 (HuobiToken) ---
Set-TextCell "E4" "0.07%"

# --- Row 5 (Cronos) ---
Set-TextCell "D5" "0.05605"

# --- Row 6 (KuCoinToken) ---
Set-TextCell "E6" "-0.09%"

# --- Row 7 (MXToken) ---
Set-TextCell "D7" "0.8189"
Set-TextCell "E7" "0.07%"

# --- Row 8 (FTXToken) ---
Set-TextCell "D8" "0.8323"
Set-TextCell "E8" "-1.01%"

# --- Row 10 (MandalaExchangeToken) ---
Set-TextCell "D10" "0.06937"
Set-TextCell "E10" "-0.16%"

# --- Row 11 (BitrueCoin) ---
Set-TextCell "D11" "0.02896"
Set-TextCell "E11" "0.92%"

# --- Row 12 (BitMartToken) ---
Set-TextCell "D12" "0.09383"
Set-TextCell "E12" "-0.01%"

# --- Row 13 (BitForexToken) ---
Set-TextCell "D13" "0.001522"
Set-TextCell "E13" "0.25%"

# --- Row 14 (One) ---
Set-TextCell "D14" "0.0006004"
Set-TextCell "E14" "-93.81%"

# --- Row 15 (TigerCash) ---
Set-TextCell "D15" "0.006199"
Set-TextCell "E15" "-0.85%"

# --- Row 16 (LEO) ---
Set-TextCell "D16" "3.655"
Set-TextCell "E16" "3.60%"

# --- Row 17 (GateToken) ---
Set-TextCell "D17" "3.025"
Set-TextCell "E17" "0.29%"

# --- Row 18 (BTSEToken) ---
Set-TextCell "D18" "2.300"
Set-TextCell "E18" "13.75%"

# --- Row 20 (LiechtensteinCryptoassetsExchange) ---
Set-TextCell "D20" "0.03088"
Set-TextCell "E20" "-3.80%"

# --- Row 21 (ProBitToken) ---
Set-TextCell "E21" "-1.48%"

# --- Row 22 (MCDex) ---
Set-TextCell "D22" "3.746"
Set-TextCell "E22" "0.10%"

# --- Row 23 (CoinExToken) ---
Set-TextCell "D23" "0.04598"
Set-TextCell "E23" "-2.32%"

# --- Row 24 (ZBToken) ---
Set-TextCell "D24" "0.1343"
Set-TextCell "E24" "-2.34%"

# --- Row 25 (BitKan) ---
Set-TextCell "D25" "0.001227"
Set-TextCell "E25" "-1.73%"

# --- Row 26 (HotbitToken) ---
Set-TextCell "D26" "0.004490"
Set-TextCell "E26" "-2.61%"

# --- Row 27 (NitroEx) ---
Set-TextCell "D27" "0.00009603"
Set-TextCell "E27" "-1.01%"

# --- Row 28 (UpBots) ---
Set-TextCell "E28" "0.73%"

# --- Row 40 (IDEX) ---
Set-TextCell "D40" "0.03640"
Set-TextCell "E40" "-0.39%"

# --- Row 41 : was BKEXToken, now KickToken ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D41" "0.006182"
Set-TextCell "E41" "0.19%"

# --- Row 42 : was KickToken, now BKEXToken ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1051"
Set-TextCell "E42" "-0.03%"

# --- Row 43 (CEJI) ---
Set-TextCell "D43" "0.002401"
Set-TextCell "E43" "-4.65%"

# --- Row 44 (LocalTraders) ---
Set-TextCell "D44" "0.008103"
Set-TextCell "E44" "5.98%"

# --- Row 45 (CoinLion) ---
Set-TextCell "D45" "0.00005355"
Set-TextCell "E45" "0.81%"

# --- Row 47 (CoinbaseStockToken) ---
Set-TextCell "D47" "0.1401"
Set-TextCell "E47" "4.93%"

# --- Row 48 (BOLO) ---
Set-TextCell "E48" "15.81%"

# --- Row 49 (CryptobidCoin) ---
Set-TextCell "D49" "0.00002101"

# --- Row 50 (SpecialPowerGold) ---
Set-TextCell "D50" "0.0002001"
